# Generate Report for Handoff
# The "1aa451c0-..." file has finished translation and is now ready for handoff,
# and a new handoff report was generated for all files that are ready for handoff
# (updating their "Latest Handoff Datetime" to the time the report was generated).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is "1aa451c0-...md", columns B (zh-cn) / C (de-de) ---
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet ---
# Row 3 ("1aa451c0-...md") status moves from "In Translation" to "Ready for handoff"
$zhcn.Range("B3").Value = "Ready for handoff"

# New handoff report timestamp applied to the "Latest Handoff Datetime" column (D)
# for every row that is ready for handoff (rows 3-7).
$zhcnHandoffDatetime = "2016-02-29 13:54:15"
$zhcn.Range("D3").Value = $zhcnHandoffDatetime
$zhcn.Range("D4").Value = $zhcnHandoffDatetime
$zhcn.Range("D5").Value = $zhcnHandoffDatetime
$zhcn.Range("D6").Value = $zhcnHandoffDatetime
$zhcn.Range("D7").Value = $zhcnHandoffDatetime

# --- de-de sheet ---
$dede.Range("B3").Value = "Ready for handoff"

$dedeHandoffDatetime = "2016-02-29 13:54:25"
$dede.Range("D3").Value = $dedeHandoffDatetime
$dede.Range("D4").Value = $dedeHandoffDatetime
$dede.Range("D5").Value = $dedeHandoffDatetime
$dede.Range("D6").Value = $dedeHandoffDatetime
$dede.Range("D7").Value = $dedeHandoffDatetime
